$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "29.311.11"
$ws.Range("D3").Value2 = "1.859.41"
$ws.Range("E3").Value2 = "  -0.01%  "
$ws.Range("E4").Value2 = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "0.7016"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value2 = "  -0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "237.47"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value2 = "  -0.09%  "
$ws.Range("E7").Value2 = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.07794"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value2 = "  -4.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.3047"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value2 = "  +0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "24.76"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value2 = "  +6.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.08139"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value2 = "  -0.49%  "
$ws.Range("D12").Value2 = "1.865.67"
$ws.Range("E12").Value2 = "  +0.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "5.208"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value2 = "  +0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.7128"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value2 = "  -0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "89.10"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value2 = "  -0.22%  "
$ws.Range("D16").Value2 = "29.282.74"
$ws.Range("E16").Value2 = "  +0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "242.97"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value2 = "  +2.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "5.769"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value2 = "  -0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "0.000007768"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value2 = "  -1.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "13.17"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value2 = "  -1.60%  "
$ws.Range("B21").Value2 = "Dai"
$ws.Range("C21").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "0.9997"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value2 = "  +0.04%  "
$ws.Range("B22").Value2 = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value2 = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value2 = "2.094.56"
$ws.Range("E22").Value2 = "  -0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "1.000"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value2 = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "7.514"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value2 = "  +0.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "162.14"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value2 = "  +0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "8.848"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value2 = "  -1.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "0.1435"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value2 = "  -1.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "18.02"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value2 = "  -0.35%  "
$ws.Range("E29").Value2 = "  -4.18%  "
$ws.Range("E30").Value2 = "  -4.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "1.472"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value2 = "  -0.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "4.287"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value2 = "  -2.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "4.024"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value2 = "  -0.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.05157"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value2 = "  -1.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.179"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value2 = "  +0.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.7057"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value2 = "  -0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.9944"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value2 = "  -0.62%  "
$ws.Range("E38").Value2 = "  +0.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.01843"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value2 = "  -0.30%  "
$ws.Range("E40").Value2 = "  -1.20%  "
$ws.Range("D41").Value2 = "1.170.43"
$ws.Range("E41").Value2 = "  +2.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.9126"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value2 = "  -1.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "5.985"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value2 = "  +0.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "71.11"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value2 = "  +0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.4238"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value2 = "  -0.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "1.000"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value2 = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "101.24"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value2 = "  -1.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.5346"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value2 = "  -1.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "1.741"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value2 = "  -2.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "9.138"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value2 = "  -0.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "6.933"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value2 = "  -0.28%  "
